$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VerifyDeleteOffer")
$ws.Activate()
$ws.Range("A27:H30").Insert(-4162) # xlShiftDown
$ws.Range("A31:H31").Copy()
$ws.Range("A27:H30").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A27").Value = "TS1_Regr_14"
$ws.Range("B27").Value = "TS_18"
$ws.Range("C27").Value = "Deselect Microsite"
$ws.Range("D27").Value = "chkInCall"
$ws.Range("E27").Value = "checkboxUncheck"

$ws.Range("A28").Value = "TS1_Regr_14"
$ws.Range("B28").Value = "TS_19"
$ws.Range("C28").Value = "Deselect Microsite"
$wsRef = $wb.Worksheets.Item("VerifyCSVForExistingVersion")
$wsRef.Range("D30").Copy()
$ws.Range("D28").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("D28").Value = "chkWeb"
$ws.Range("E28").Value = "checkboxUncheck"

$ws.Range("A29").Value = "TS1_Regr_14"
$ws.Range("B29").Value = "TS_20"
$ws.Range("C29").Value = "Deselect Microsite"
$ws.Range("D29").Value = "chkDirectMail"
$ws.Range("E29").Value = "checkboxUncheck"

$ws.Range("A30").Value = "TS1_Regr_14"
$ws.Range("B30").Value = "TS_21"
$ws.Range("C30").Value = "Deselect Microsite"
$ws.Range("D30").Value = "chkCallCentre"
$ws.Range("E30").Value = "checkboxUncheck"

$ws.Range("F35").Value = "free_mobile"

$excel.ActiveWindow.ScrollRow = 31
$ws.Range("F45").Select()

$wsNew = $wb.Worksheets.Item("VerifyCSVForNewVersion")
$wsNew.Activate()
$wsNew.Range("D27").Select()

$ws.Activate()

Write-Host "done"
